# BIS-768: Fixed XLS export test files
# Adds a new "Unique" column (L) to the experiment-type property-type table,
# mirroring the existing "Multivalued" column (K): header in row 4, and a
# "FALSE" value copied down for each property-type data row (5-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell L4: copy K4's look (bold header style) and set new label ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L4").Value2 = "Unique"

# --- Data cells L5:L8: copy K5:K8's look (TRUE/FALSE format) and set value ---
# Leading apostrophe forces literal text "FALSE" (matching column K's stored
# type) instead of Excel auto-coercing the word into a boolean value.
$ws.Range("K5:K8").Copy()
$ws.Range("L5:L8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L5").Value2 = "'FALSE"
$ws.Range("L6").Value2 = "'FALSE"
$ws.Range("L7").Value2 = "'FALSE"
$ws.Range("L8").Value2 = "'FALSE"

$excel.CutCopyMode = 0

# --- Update selection to mirror the author's last-touched range ---
$ws.Range("L7:L8").Select() | Out-Null
